# Template_upload_data_leads-2.xlsx — "Fixing and update to digisales mobile version 2"
#
# Business-data changes:
#  - Column A (CIF) rows 2-8: bump the base id from 9020304198 to 9020304205.
#    A3:A8 are formulas chained off the previous row (=prev+1 / shared
#    formula), so updating A2 alone makes the whole block follow, same as
#    it would in Excel.
#  - Column B (Nama Cust) rows 2-8: "dedic 36".."dedic 42" -> "dedic 71".."dedic 77".
#  - Column O (Expired Date, =TODAY()+1) is a volatile formula; it
#    recalculates itself and needs no manual edit.
#
# View-state changes (best effort via the Excel object model):
#  - Zoom to 100%.
#  - Selection moves from D10 to S6.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column A: renumber the CIF id sequence -------------------------------
$ws.Range("A2").Value = 9020304205

# --- Column B: update the "dedic NN" customer-name placeholders -----------
$newNames = @("dedic 71", "dedic 72", "dedic 73", "dedic 74", "dedic 75", "dedic 76", "dedic 77")
for ($i = 0; $i -lt $newNames.Length; $i++) {
    $row = 2 + $i
    $ws.Cells.Item($row, 2).Value = $newNames[$i]
}

# --- View state: zoom + selection ------------------------------------------
$win = $excel.ActiveWindow
$win.Zoom = 100
$ws.Range("S6").Select() | Out-Null
